# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 2
# of the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-19 04:54:35"
$zhcn.Range("H2").Value = "2016-03-19 04:54:53"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-19 04:54:38"
$dede.Range("H2").Value = "2016-03-19 04:54:58"
